$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header in G1, matching the bold style used by the other headers
$ws.Range("G1").Value = "SOMAR IDEIAS"
$ws.Range("G1").Font.Bold = $true

# Add new data values
$ws.Cells.Item(4, 7).Value = 24
$ws.Cells.Item(10, 7).Value = 4

# Update selection to G2
$ws.Range("G2").Select()
